# Fix project gantt chart
# Remove the "EstimatedEffortHours" (col C) and "Progress" (col E) columns
# from the Projects table / worksheet, keeping ID, ProjectName, Dependencies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the table binding first so the column deletes below don't leave the
# table definition out of sync; we'll recreate the table afterwards.
$lo = $ws.ListObjects.Item(1)
$lo.Unlist()

# Delete the rightmost column first (Progress, column E) and then the
# EstimatedEffortHours column (column C) so that the Dependencies column
# shifts left into column C.
$ws.Columns.Item(5).Delete()
$ws.Columns.Item(3).Delete()

# Recreate the table over the remaining data (ID, ProjectName, Dependencies).
$newTable = $ws.ListObjects.Add(1, $ws.Range("A1:C3"), [System.Reflection.Missing]::Value, 1, "")
$newTable.Name = "Table1"

# Update the selection like the authored workbook.
$ws.Range("C5").Select()
